$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6491.7896
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 6835.778
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 20507.334
$ws.Range("M17").Value = -732
$ws.Range("N17").Value = -20843.334

# ALC row 26
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2913
$ws.Range("I113").Value = 1840.8
$ws.Range("J113").Value = 3449.1
$ws.Range("K113").Value = 1840.8
$ws.Range("L113").Value = 3449.1
$ws.Range("M113").Value = 1413.2
$ws.Range("N113").Value = -9957.1

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3230.84
$ws.Range("I137").Value = 1745.55
$ws.Range("J137").Value = 9172
$ws.Range("K137").Value = 5236.65
$ws.Range("L137").Value = 27516
$ws.Range("M137").Value = -2686.65

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 105790
$ws.Range("I140").Value = 70000
$ws.Range("J140").Value = 114737.5
$ws.Range("K140").Value = 70000
$ws.Range("L140").Value = 114737.5
$ws.Range("M140").Value = -64820
$ws.Range("N140").Value = -125097.5

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1296.8334
$ws.Range("I2").Value = 1256.9565
$ws.Range("J2").Value = 1427.8572
$ws.Range("K2").Value = 1256.9565
$ws.Range("L2").Value = 1427.8572
$ws.Range("M2").Value = -1143.9565

# ARM row 21
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 995
$ws.Range("I21").Value = 995
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 995
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -621
$ws.Range("N21").ClearContents()

# ARM row 27
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 43336
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 43336
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 43336
$ws.Range("N27").Value = -43704

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6226.107
$ws.Range("I32").Value = 4873.959
$ws.Range("J32").Value = 15691.143
$ws.Range("K32").Value = 4873.959
$ws.Range("L32").Value = 15691.143
$ws.Range("M32").Value = -4586.959

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3797.0698
$ws.Range("I61").Value = 3263.742
$ws.Range("J61").Value = 5174.8335
$ws.Range("K61").Value = 3263.742
$ws.Range("L61").Value = 5174.8335
$ws.Range("M61").Value = -3051.742
$ws.Range("N61").Value = -5598.8335

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6031.5415
$ws.Range("I74").Value = 1761.0454
$ws.Range("J74").Value = 53007
$ws.Range("K74").Value = 1761.0454
$ws.Range("L74").Value = 53007
$ws.Range("M74").Value = -887.0454
$ws.Range("N74").Value = -54755

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6031.5415
$ws.Range("I77").Value = 1761.0454
$ws.Range("J77").Value = 53007
$ws.Range("K77").Value = 8805.226999999999
$ws.Range("L77").Value = 265035
$ws.Range("M77").Value = -4437.226999999999
$ws.Range("N77").Value = -273771

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1296.8334
$ws.Range("I116").Value = 1256.9565
$ws.Range("J116").Value = 1427.8572
$ws.Range("K116").Value = 1256.9565
$ws.Range("L116").Value = 1427.8572
$ws.Range("M116").Value = 1037.0435

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3797.0698
$ws.Range("I136").Value = 3263.742
$ws.Range("J136").Value = 5174.8335
$ws.Range("K136").Value = 9791.226000000001
$ws.Range("L136").Value = 15524.5005
$ws.Range("M136").Value = -7241.226000000001
$ws.Range("N136").Value = -20624.5005

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1296.8334
$ws.Range("I3").Value = 1256.9565
$ws.Range("J3").Value = 1427.8572
$ws.Range("K3").Value = 1256.9565
$ws.Range("L3").Value = 1427.8572
$ws.Range("M3").Value = -1142.9565

# BSM row 34
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

# BSM row 46
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# CRP row 23
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 51994.875
$ws.Range("I23").Value = 1969.6666
$ws.Range("J23").Value = 82010
$ws.Range("K23").Value = 1969.6666
$ws.Range("L23").Value = 82010
$ws.Range("M23").Value = -1729.6666
$ws.Range("N23").Value = -82490

# CRP row 27
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 51994.875
$ws.Range("I27").Value = 1969.6666
$ws.Range("J27").Value = 82010
$ws.Range("K27").Value = 1969.6666
$ws.Range("L27").Value = 82010
$ws.Range("M27").Value = -1777.6666
$ws.Range("N27").Value = -82394

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5755.6113
$ws.Range("I31").Value = 5904.037
$ws.Range("J31").Value = 5310.3335
$ws.Range("K31").Value = 5904.037
$ws.Range("L31").Value = 5310.3335
$ws.Range("M31").Value = -5609.037
$ws.Range("N31").Value = -5900.3335

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5755.6113
$ws.Range("I34").Value = 5904.037
$ws.Range("J34").Value = 5310.3335
$ws.Range("K34").Value = 5904.037
$ws.Range("L34").Value = 5310.3335
$ws.Range("M34").Value = -5702.037
$ws.Range("N34").Value = -5714.3335

# CRP row 130
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H130").Value = 60468
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 60468
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 60468
$ws.Range("N130").Value = -70508

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1635.9762
$ws.Range("I132").Value = 1322.8334
$ws.Range("J132").Value = 2418.8333
$ws.Range("K132").Value = 3968.5002
$ws.Range("L132").Value = 7256.499899999999
$ws.Range("M132").Value = -1438.5002
$ws.Range("N132").Value = -12316.4999

# CRP row 138
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 44833.332
$ws.Range("I138").Value = 30000
$ws.Range("J138").Value = 52250
$ws.Range("K138").Value = 30000
$ws.Range("L138").Value = 52250
$ws.Range("M138").Value = -24860
$ws.Range("N138").Value = -62530

# CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# CUL row 19
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 3000.25
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 4000
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 171
$ws.Range("N19").Value = -12348

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 778.87
$ws.Range("I113").Value = 784.5349
$ws.Range("J113").Value = 744.0714
$ws.Range("K113").Value = 2353.6047
$ws.Range("L113").Value = 2232.2142
$ws.Range("M113").Value = -183.6046999999999
$ws.Range("N113").Value = -6572.2142

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3381.4614
$ws.Range("I141").Value = 1518.1666
$ws.Range("J141").Value = 4978.5713
$ws.Range("K141").Value = 4554.4998
$ws.Range("L141").Value = 14935.7139
$ws.Range("M141").Value = 625.5002000000004
$ws.Range("N141").Value = -25295.7139

# LTW row 14
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 10000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 10000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -10344

# LTW row 21
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7855.95
$ws.Range("I61").Value = 10894.4
$ws.Range("J61").Value = 4817.5
$ws.Range("K61").Value = 10894.4
$ws.Range("L61").Value = 4817.5
$ws.Range("M61").Value = -10692.4
$ws.Range("N61").Value = -5221.5

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 7855.95
$ws.Range("I113").Value = 10894.4
$ws.Range("J113").Value = 4817.5
$ws.Range("K113").Value = 10894.4
$ws.Range("L113").Value = 4817.5
$ws.Range("M113").Value = -8724.4
$ws.Range("N113").Value = -9157.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4636.314
$ws.Range("I122").Value = 4385.05
$ws.Range("J122").Value = 5550
$ws.Range("K122").Value = 13155.15
$ws.Range("L122").Value = 16650
$ws.Range("M122").Value = -10705.15
$ws.Range("N122").Value = -21550

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5998.472
$ws.Range("I136").Value = 4346.731
$ws.Range("J136").Value = 10293
$ws.Range("K136").Value = 13040.193
$ws.Range("L136").Value = 30879
$ws.Range("M136").Value = -10490.193
$ws.Range("N136").Value = -35979

# WVR row 26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 7620
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 7620
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 7620
$ws.Range("N26").Value = -8206

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4548.364
$ws.Range("I136").Value = 1730.0286
$ws.Range("J136").Value = 9480.450000000001
$ws.Range("K136").Value = 5190.085800000001
$ws.Range("L136").Value = 28441.35
$ws.Range("M136").Value = -2640.085800000001
$ws.Range("N136").Value = -33541.35000000001
